$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.627098560333252
$ws.Range("B1").Value = 2.881130218505859
$ws.Range("C1").Value = 1.734490752220154
$ws.Range("D1").Value = 1.408208847045898
$ws.Range("E1").Value = 1.313300967216492
